# Update cryptocurrency price/volume data per the upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "34.492.17"
$ws.Cells.Item(2, 5).Value = "  +0.31%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.808.68"
$ws.Cells.Item(3, 5).Value = "  +0.45%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "225.09"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.601"
$ws.Cells.Item(6, 5).Value = "  +5.05%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "38.43"
$ws.Cells.Item(8, 5).Value = "  +6.16%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -3.95%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0670"
$ws.Cells.Item(10, 5).Value = "  -2.88%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0981"
$ws.Cells.Item(11, 5).Value = "  +1.73%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "2.069.06"
$ws.Cells.Item(12, 5).Value = "  +0.50%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "11.16"
$ws.Cells.Item(13, 5).Value = "  -3.27%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "1.804.51"
$ws.Cells.Item(14, 5).Value = "  +0.28%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  -1.82%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "34.461.49"
$ws.Cells.Item(16, 5).Value = "  +0.36%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "4.38"
$ws.Cells.Item(17, 5).Value = "  -2.43%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "68.23"
$ws.Cells.Item(18, 5).Value = "  -0.93%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "241.30"
$ws.Cells.Item(19, 5).Value = "  -1.28%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "0.0₃0769"
$ws.Cells.Item(20, 5).Value = "  -2.62%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "11.18"
$ws.Cells.Item(21, 5).Value = "  -3.48%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  -0.09%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -1.53%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  +1.50%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "170.66"
$ws.Cells.Item(25, 5).Value = "  -1.13%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "7.69"
$ws.Cells.Item(26, 5).Value = "  -3.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "17.52"
$ws.Cells.Item(27, 5).Value = "  +4.32%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.121"
$ws.Cells.Item(28, 5).Value = "  +3.27%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -0.04%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -0.85%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "3.78"
$ws.Cells.Item(31, 5).Value = "  -1.26%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "Hedera"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.0515"
$ws.Cells.Item(32, 5).Value = "  -2.29%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "3.85"
$ws.Cells.Item(33, 5).Value = "  -4.07%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +1.22%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "ImmutableX"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.640"
$ws.Cells.Item(35, 5).Value = "  -4.72%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "Maker"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(36, 4).Value = "1.308.40"
$ws.Cells.Item(36, 5).Value = "  -6.25%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  -0.52%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -1.56%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -4.60%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "82.53"

# Row 41
$ws.Cells.Item(41, 2).Value = "HuobiToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.44"
$ws.Cells.Item(41, 5).Value = "  +1.03%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "WEMIXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.22"
$ws.Cells.Item(42, 5).Value = "  +3.28%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -0.07%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -0.84%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "14.06"
$ws.Cells.Item(45, 5).Value = "  +5.54%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  +2.42%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "1.970.77"
$ws.Cells.Item(47, 5).Value = "  +0.55%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "5.81"
$ws.Cells.Item(48, 5).Value = "  -3.26%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -0.11%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "102.97"
$ws.Cells.Item(50, 5).Value = "  -1.11%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  -5.35%  "

